$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.490.34'
$ws.Range('E2').Value = '  +0.81%  '
$ws.Range('D3').Value = '2.606.85'
$ws.Range('E3').Value = '  +1.03%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''537.58'
$ws.Range('E5').Value = '  +3.33%  '
$ws.Range('D6').Value = '''141.51'
$ws.Range('E6').Value = '  +2.28%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '''0.567'
$ws.Range('E8').Value = '  +0.82%  '
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('E10').Value = '  +1.70%  '
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').Value = '3.066.60'
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('D14').Value = '59.407.45'
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').Value = '2.624.65'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').Value = '''341.39'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('D21').Value = '''6.37'
$ws.Range('E21').Value = '  -1.92%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '''67.46'
$ws.Range('E23').Value = '  +2.28%  '
$ws.Range('E24').Value = '  +1.81%  '
$ws.Range('E25').Value = '  -1.37%  '
$ws.Range('D26').Value = '''0.994'
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('E27').Value = '  +3.54%  '
$ws.Range('D28').Value = '0.0₃0746'
$ws.Range('E28').Value = '  +3.79%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  +6.27%  '
$ws.Range('E31').Value = '  -1.65%  '
$ws.Range('E32').Value = '  +1.38%  '
$ws.Range('D33').Value = '''149.55'
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').Value = '''1.46'
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('D37').Value = '''0.834'
$ws.Range('E37').Value = '  +3.95%  '
$ws.Range('D38').Value = '''0.827'
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('E39').Value = '  +1.23%  '
$ws.Range('D40').Value = '''0.999'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').Value = '''273.59'
$ws.Range('E41').Value = '  +1.59%  '
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('E43').Value = '  -0.32%  '
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('E45').Value = '  +1.48%  '
$ws.Range('D46').Value = '1.950.52'
$ws.Range('E46').Value = '  -0.47%  '
$ws.Range('D47').Value = '''18.52'
$ws.Range('E47').Value = '  +3.89%  '
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('D50').Value = '''110.89'
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('D51').Value = '''4.74'
$ws.Range('E51').Value = '  +0.45%  '
